$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 45, shifting existing rows 45-49 down to 46-50
$ws.Rows("45:45").Insert() | Out-Null

# Populate the new row with course name and rating value
$ws.Range("A45").Value = "Introduction to Statistics in Python"
$ws.Range("I45").Value = 3

# Update the active selection to A46, matching the recorded cursor position
$ws.Range("A46").Select() | Out-Null
